# Updated queries for C3DC first half testcases.
# The SQL queries stored in the "StatQuery" (C2) and "TabQuery" (B2:B7) cells
# joined tables using the generic "id" column. Update them to join on the
# fully-qualified study_id / participant_id columns instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellAddresses = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellAddresses) {
    $cell = $ws.Range($addr)
    $text = $cell.Text

    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $cell.Value = $text
}
